# Contemplar credenciales no validas
# The "locked_out_user" row is removed, and the remaining rows are
# reorganized: problem_user moves to row 2, performance_glitch_user moves
# to row 3. The hyperlink that pointed at the old "problem_user" row
# (correo03@gmail.com) is removed; the other two hyperlinks (correo01,
# correo02) are kept on A1/A2. The sheet now only spans A1:B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that used to sit on A3 (correo03@gmail.com). This
# runtime's Hyperlinks.Delete() only operates on the whole sheet
# collection, so drop them all and re-create the two that must survive
# (A1/correo01, A2/correo02), then strip the auto-applied "Hyperlink"
# cell style back off so the cells keep their original (default)
# formatting.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:correo01@gmail.com", "", "", "correo01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:correo02@gmail.com", "", "", "correo02@gmail.com")
$ws.Range("A1").Style = "Normal"
$ws.Range("A2").Style = "Normal"

# Rewrite the data area to reflect the new 3-row layout.
$ws.Range("A1").Value = "standard_user"
$ws.Range("B1").Value = "secret_sauce"

$ws.Range("A2").Value = "problem_user"
$ws.Range("B2").Value = "secret_sauce"

$ws.Range("A3").Value = "performance_glitch_user"
$ws.Range("B3").Value = "secret_sauce"

# Row 4 is no longer part of the table; clear its old contents entirely.
$ws.Range("A4:B4").Clear()

# Update the saved selection to match the authored state.
$ws.Range("A7").Select()
